$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 62 (inclusive) in column C currently hold 7310 and must
# become 7293, matching the already-7293 values found in rows 63-252.
$ws.Range("C2:C62").Value = 7293
